# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 07:22"

# 2. Row 30 (Pakistan) - update case numbers
$ws.Range("B30").Value = 14079
$ws.Range("C30").Value = 164
$ws.Range("D30").Value = 3233
$ws.Range("E30").Value = 10545
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = 301

# 3. Row 61 (Tailandia) - update case numbers
$ws.Range("B61").Value = 2938
$ws.Range("C61").Value = 7
$ws.Range("D61").Value = 2652
$ws.Range("E61").Value = 232
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 54

# 4. Rows 83/84 - Bulgaria and Eslovaquia swap order (Bulgaria moves up with
#    fresh data, Eslovaquia moves down keeping its former data)
$ws.Range("A83").Value = "Bulgaria"
$ws.Range("B83").Value = 1387
$ws.Range("C83").Value = 24
$ws.Range("D83").Value = 222
$ws.Range("E83").Value = 1107
$ws.Range("F83").Value = 41
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 58

$ws.Range("A84").Value = "Eslovaquia"
$ws.Range("B84").Value = 1381
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 403
$ws.Range("E84").Value = 960
$ws.Range("F84").Value = 7
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 18

# 5. Row 152 (Zambia) - update case numbers
$ws.Range("C152").Value = 0

$wb.Save()
